$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("AF", 20, "ADS", 182983),
    @("AN", 22, "ADS", 102547),
    @("CL", 16, "ADS", 113274),
    @("DM", 12, "ADS", 102641),
    @("IA", 23, "ADS", 250898),
    @("JA", 24, "ADS", 264134),
    @("JE", 21, "ADS", 106942),
    @("KA", 14, "ADS", 201651),
    @("KV", 16, "ADS", 86234),
    @("LK", 16, "ADS", 220149),
    @("LM", 12, "ADS", 114389),
    @("LV", 16, "ADS", 108247),
    @("MA", 23, "ADS", 281751),
    @("ME", 12, "ADS", 155431),
    @("MPS", 14, "ADS", 120351),
    @("MS", 23, "ADS", 155818),
    @("MV", 16, "ADS", 132375),
    @("PV", 11, "ADS", 220867),
    @("RV", 11, "ADS", 120169),
    @("SA", 15, "ADS", 108956),
    @("SAF", 19, "ADS", 235666),
    @("SC", 19, "ADS", 145503),
    @("SS", 14, "ADS", 170955),
    @("TL", 11, "ADS", 230585),
    @("TM", 24, "ADS", 201805)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row++
}
